$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, pushing existing rows 17-29 down to 18-30.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with its values.
$ws.Cells.Item(17, 1).Value = 2
$ws.Cells.Item(17, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44447
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100112022
$ws.Cells.Item(17, 7).Value = "Arveja Verde"
$ws.Cells.Item(17, 8).Value = "Perfection"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 600
$ws.Cells.Item(17, 11).Value = 28000
$ws.Cells.Item(17, 12).Value = 30000
$ws.Cells.Item(17, 13).Value = 29000
$ws.Cells.Item(17, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 16).Value = 1160
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Match the date-column formatting used by the other rows in the table
# (only column D carries an explicit style, the same as the row below it).
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat
